$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# --- Recommandations sheet updates ---
$ws.Range("C2").Value = 10
$ws.Range("D2").Value = 4070.22
$ws.Range("E2").Value = 105.01
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 3790
$ws.Range("E3").Value = 860
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 3270.54
$ws.Range("E4").Value = 658.16
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 3260
$ws.Range("E5").Value = 655
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 2910
$ws.Range("E6").Value = 590
$ws.Range("C7").Value = 5
$ws.Range("D7").Value = 2890
$ws.Range("E7").Value = 595
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 2835
$ws.Range("C9").Value = 5
$ws.Range("D9").Value = 2585
$ws.Range("E9").Value = 525
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 1796.12
$ws.Range("E10").Value = 362.92
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 1739.1
$ws.Range("C12").Value = 5
$ws.Range("D12").Value = 1563.78
$ws.Range("E12").Value = 316.23
$ws.Range("C13").Value = 5
$ws.Range("D13").Value = 1254.68
$ws.Range("E13").Value = 261.41
$ws.Range("C14").Value = 5
$ws.Range("D14").Value = 1037.85
$ws.Range("E14").Value = 215.21
$ws.Range("C15").Value = 5
$ws.Range("D15").Value = 931.52
$ws.Range("E15").Value = 189.01
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 684.0700000000001
$ws.Range("E16").Value = 138.9
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 647.78
$ws.Range("E17").Value = 129.27
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 611.6799999999999
$ws.Range("E18").Value = 121.93
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 601.16
$ws.Range("E19").Value = 119.83
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 530.92
$ws.Range("E20").Value = 105.5
$ws.Range("C21").Value = 5
$ws.Range("D21").Value = 526.14
$ws.Range("E21").Value = 107.45
$ws.Range("C22").Value = 5
$ws.Range("D22").Value = 461.44
$ws.Range("E22").Value = 92.56999999999999
$ws.Range("B23").Value = 5
$ws.Range("D23").Value = 37.36
$ws.Range("E23").Value = 7.49
$ws.Range("B24").Value = 3
$ws.Range("D24").Value = 22.07
$ws.Range("E24").Value = 7.5
$ws.Range("F24").Value = "🟢 Achat"
$ws.Range("G24").Value = "✅ Renforcer"
$ws.Range("A27").Value = "SICABLE CI (CABC)"
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 4.17
$ws.Range("E27").Value = 4.17
$ws.Range("G27").Value = "➖ Neutre"
$ws.Range("A28").Value = "ORAGROUP TOGO (ORGT)"
$ws.Range("D28").Value = 4.09
$ws.Range("E28").Value = -1.74
$ws.Range("A29").Value = "TRACTAFRIC MOTORS CI (PRSC)"
$ws.Range("B29").Value = 1
$ws.Range("D29").Value = 3.98
$ws.Range("E29").Value = -3.42
$ws.Range("A30").Value = "SOCIETE GENERALE COTE D'IVOIRE (SGBC)"
$ws.Range("D30").Value = 3.58
$ws.Range("E30").Value = 7.04
$ws.Range("A31").Value = "SOCIETE IVOIRIENNE DE BANQUE  (SIBC)"
$ws.Range("B31").Value = 2
$ws.Range("C31").Value = 1
$ws.Range("D31").Value = 2.88
$ws.Range("E31").Value = -3.03
$ws.Range("G31").Value = "👀 À surveiller"
$ws.Range("A32").Value = "BANK OF AFRICA SENEGAL (BOAS)"
$ws.Range("B32").Value = 1
$ws.Range("D32").Value = 2.43
$ws.Range("E32").Value = -2.22
$ws.Range("A33").Value = "SOGB CI (SOGC)"
$ws.Range("C33").Value = 0
$ws.Range("D33").Value = 1.8
$ws.Range("E33").Value = 1.8
$ws.Range("G33").Value = "➖ Neutre"
$ws.Range("A34").Value = "SETAO CI (STAC)"
$ws.Range("D34").Value = 1.77
$ws.Range("E34").Value = 3.48
$ws.Range("A35").Value = "BERNABE CI (BNBC)"
$ws.Range("B35").Value = 2
$ws.Range("C35").Value = 1
$ws.Range("D35").Value = 1.18
$ws.Range("E35").Value = 3.59
$ws.Range("G35").Value = "👀 À surveiller"
$ws.Range("A36").Value = "CIE CI (CIEC)"
$ws.Range("D36").Value = 0.51
$ws.Range("E36").Value = -5.06
$ws.Range("A37").Value = "BANK OF AFRICA BN (BOAB)"
$ws.Range("D37").Value = 0.42
$ws.Range("E37").Value = 2.86
$ws.Range("A38").Value = "TOTAL"
$ws.Range("C38").Value = 5
$ws.Range("D38").Value = 0
$ws.Range("E38").Value = 0
$ws.Range("A39").Value = "UNIWAX CI (UNXC)"
$ws.Range("B39").Value = 1
$ws.Range("D39").Value = -0.71
$ws.Range("E39").Value = 3.6
$ws.Range("G39").Value = "👀 À surveiller"
$ws.Range("A40").Value = "SOLIBRA CI (SLBC)"
$ws.Range("D40").Value = -0.8100000000000001
$ws.Range("E40").Value = -0.8100000000000001
$ws.Range("A41").Value = "NEI-CEDA CI (NEIC)"
$ws.Range("D41").Value = -0.84
$ws.Range("E41").Value = -0.84
$ws.Range("A42").Value = "FILTISAC CI (FTSC)"
$ws.Range("B42").Value = 1
$ws.Range("C42").Value = 2
$ws.Range("D42").Value = -1.41
$ws.Range("E42").Value = -0.8
$ws.Range("G42").Value = "👀 À surveiller"
$ws.Range("A43").Value = "TOTALENERGIES MARKETING CI (TTLC)"
$ws.Range("D43").Value = -1.6
$ws.Range("E43").Value = -1.6
$ws.Range("A44").Value = "VIVO ENERGY CI (SHEC)"
$ws.Range("D44").Value = -2
$ws.Range("E44").Value = -2
$ws.Range("A45").Value = "ONATEL BF (ONTBF)"
$ws.Range("D45").Value = -2.08
$ws.Range("E45").Value = -2.08
$ws.Range("A46").Value = "SITAB CI (STBC)"
$ws.Range("D46").Value = -2.5
$ws.Range("E46").Value = -2.5
$ws.Range("A47").Value = "SICOR CI (SICC)"
$ws.Range("C47").Value = 1
$ws.Range("D47").Value = -2.73
$ws.Range("E47").Value = -2.73
$ws.Range("A48").Value = "CFAO MOTORS CI (CFAC)"
$ws.Range("B48").Value = 0
$ws.Range("C48").Value = 1
$ws.Range("D48").Value = -3.65
$ws.Range("E48").Value = -3.65
$ws.Range("F48").Value = "🟡 Observer"
$ws.Range("G48").Value = "➖ Neutre"
$ws.Range("A49").Value = "SONATEL SN (SNTS)"
$ws.Range("B49").Value = 0
$ws.Range("C49").Value = 1
$ws.Range("D49").Value = -3.81
$ws.Range("E49").Value = -3.81
$ws.Range("F49").Value = "🟡 Observer"
$ws.Range("G49").Value = "➖ Neutre"
$ws.Range("A50").Value = "SUCRIVOIRE (SCRC)"
$ws.Range("B50").Value = 0
$ws.Range("C50").Value = 1
$ws.Range("D50").Value = -3.85
$ws.Range("E50").Value = -3.85
$ws.Range("F50").Value = "🟡 Observer"
$ws.Range("G50").Value = "➖ Neutre"
$ws.Range("A51").Value = "BANK OF AFRICA NG (BOAN)"
$ws.Range("B51").Value = 0
$ws.Range("C51").Value = 3
$ws.Range("D51").Value = -9.24
$ws.Range("E51").Value = -2.19
$ws.Range("F51").Value = "🔴 Vente"
$ws.Range("G51").Value = "⚠️ Risque de décrochage"

# --- Top_YTD sheet updates ---
$ws2.Range("B2").Value = 126993511.08
$ws2.Range("B3").Value = 4585167.44
$ws2.Range("B4").Value = 2438343.51
$ws2.Range("B5").Value = 2403418.91
$ws2.Range("B6").Value = 1475135.18
$ws2.Range("B7").Value = 1430169.86
$ws2.Range("B8").Value = 1318979.94
$ws2.Range("B9").Value = 893515.72
$ws2.Range("B10").Value = 204120.96
$ws2.Range("B11").Value = 179966.07
